$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"5.075712372040022e-17"
$ws.Range("D3").Value = [double]"6.175624517019239e-17"
$ws.Range("D4").Value = [double]"0.1389898938418274"
$ws.Range("D5").Value = [double]"0.131042139217673"
$ws.Range("D6").Value = [double]"0.09787673297163955"
$ws.Range("D7").Value = [double]"0.6320912339688598"
$ws.Range("D8").Value = [double]"6.24064294285224e-17"
